# The document had a couple of sentences split across three runs each
# (likely from earlier edits). This collapses each of those sentences
# back into a single contiguous run of text, with no wording change.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "The result of the crash ended up with a disk that appeared to be total empty. It was visible in the directory listing as an empty drive.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The result of the crash ended up with a disk that appeared to be total empty. It was visible in the directory listing as an empty drive.",
    2
)

$d.Content.Find.Execute(
    "This lists the available usb drives. Any NVMe system drives, which may show up as usb drives are excluded. Selecting one of the drives will show a number of additional tabs and open the Boot Parameter Table tab.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "This lists the available usb drives. Any NVMe system drives, which may show up as usb drives are excluded. Selecting one of the drives will show a number of additional tabs and open the Boot Parameter Table tab.",
    2
)
